$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K4 value (65 -> 70); dependent formulas K5 and L5 will recalc automatically.
$ws.Range("K4").Value = 70

# Update the selected cell shown in the sheet view.
$ws.Range("I10").Select()

# Add new row 12 data.
# Copy the date style (numFmtId 14) from an existing date cell (B11) onto B12
# so it reuses the existing style index instead of creating a new one.
$ws.Range("B11").Copy($ws.Range("B12"))
$ws.Range("B12").Value = 44446
$ws.Range("C12").Value = "late QS"
$ws.Range("D12").Value = "Day58"
$ws.Range("E12").Value = "Minh"
$ws.Range("F12").Value = 5
